$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 61759
$ws.Range("G5").Value = 60918
$ws.Range("G6").Value = 54164
$ws.Range("G7").Value = 52958
$ws.Range("G8").Value = 52913
$ws.Range("G9").Value = 56672
$ws.Range("G10").Value = 56947
$ws.Range("G11").Value = 58195
$ws.Range("G12").Value = 57720
$ws.Range("G13").Value = 58964
$ws.Range("G14").Value = 20797
$ws.Range("G15").Value = 20644
$ws.Range("G16").Value = 20445
$ws.Range("G17").Value = 20424
$ws.Range("G18").Value = 20371
$ws.Range("G19").Value = 57
$ws.Range("G20").Value = 25
$ws.Range("G21").Value = 22
$ws.Range("G22").Value = 22
$ws.Range("G23").Value = 26
$ws.Range("G24").Value = 317
$ws.Range("G25").Value = 65
$ws.Range("G26").Value = 46
$ws.Range("G27").Value = 58
$ws.Range("G28").Value = 60
$ws.Range("G29").Value = 78
$ws.Range("G30").Value = 224
$ws.Range("G31").Value = 130
$ws.Range("G32").Value = 78
$ws.Range("G33").Value = 115
$ws.Range("G34").Value = 17778
$ws.Range("G35").Value = 19389
$ws.Range("G36").Value = 23968
$ws.Range("G37").Value = 23653
$ws.Range("G38").Value = 24275
$ws.Range("G39").Value = 7
$ws.Range("G40").Value = 4
$ws.Range("G41").Value = 2
$ws.Range("G42").Value = 2
$ws.Range("G43").Value = 2
$ws.Range("G44").Value = 3
$ws.Range("G45").Value = 2
$ws.Range("G46").Value = 1
$ws.Range("G47").Value = 1
$ws.Range("G48").Value = 1
$ws.Range("G49").Value = 33
$ws.Range("G50").Value = 11
$ws.Range("G51").Value = 20
$ws.Range("G52").Value = 19
$ws.Range("G53").Value = 18
$ws.Range("G54").Value = 10
$ws.Range("G55").Value = 14
$ws.Range("G56").Value = 4
$ws.Range("G57").Value = 3
$ws.Range("G58").Value = 3
$ws.Range("G59").Value = 11
$ws.Range("G60").Value = 17
$ws.Range("G61").Value = 6
$ws.Range("G62").Value = 3
$ws.Range("G63").Value = 3

$ws.Range("G64").Select()
